$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '332.91'
Set-TextValue $ws.Range('E2') '1.12%'
Set-TextValue $ws.Range('D3') '45.60'
Set-TextValue $ws.Range('E3') '3.14%'
Set-TextValue $ws.Range('D4') '5.477'
Set-TextValue $ws.Range('E4') '-0.54%'
Set-TextValue $ws.Range('D5') '0.08531'
Set-TextValue $ws.Range('E5') '5.58%'
Set-TextValue $ws.Range('D6') '2.079'
Set-TextValue $ws.Range('E6') '1.80%'
$ws.Range('B7').Value = 'GateToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue $ws.Range('D7') '4.447'
Set-TextValue $ws.Range('E7') '0.94%'
$ws.Range('B8').Value = 'MXToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws.Range('D8') '0.9893'
Set-TextValue $ws.Range('E8') '3.70%'
$ws.Range('B9').Value = 'BTSEToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue $ws.Range('D9') '2.534'
Set-TextValue $ws.Range('E9') '-2.42%'
$ws.Range('B10').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C10').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue $ws.Range('D10') '0.1163'
Set-TextValue $ws.Range('E10') '1.24%'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue $ws.Range('D11') '0.1920'
Set-TextValue $ws.Range('E11') '1.92%'
$ws.Range('B12').Value = 'MCDex'
$ws.Range('C12').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
Set-TextValue $ws.Range('D12') '9.455'
Set-TextValue $ws.Range('E12') '-7.05%'
$ws.Range('B13').Value = 'MandalaExchangeToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue $ws.Range('D13') '0.09743'
Set-TextValue $ws.Range('E13') '-1.82%'
$ws.Range('B14').Value = 'BitrueCoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue $ws.Range('D14') '0.04724'
Set-TextValue $ws.Range('E14') '-2.09%'
$ws.Range('B15').Value = 'BitMartToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue $ws.Range('D15') '0.1059'
Set-TextValue $ws.Range('E15') '-0.49%'
$ws.Range('B16').Value = 'BitForexToken'
$ws.Range('C16').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue $ws.Range('D16') '0.001294'
Set-TextValue $ws.Range('E16') '0.89%'
$ws.Range('B17').Value = 'TigerCash'
$ws.Range('C17').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue $ws.Range('D17') '0.006006'
Set-TextValue $ws.Range('E17') '2.90%'
$ws.Range('B18').Value = 'HotbitToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
Set-TextValue $ws.Range('D18') '0.004467'
Set-TextValue $ws.Range('E18') '2.69%'
$ws.Range('B19').Value = 'LEO'
$ws.Range('C19').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws.Range('D19') '3.381'
Set-TextValue $ws.Range('E19') '0.17%'
$ws.Range('B20').Value = 'BitpandaEcosystemToken'
$ws.Range('C20').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
Set-TextValue $ws.Range('D20') '0.3355'
Set-TextValue $ws.Range('E20') '-1.51%'
$ws.Range('B21').Value = 'ProBitToken'
$ws.Range('C21').Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
Set-TextValue $ws.Range('D21') '0.1374'
Set-TextValue $ws.Range('E21') '-1.86%'
$ws.Range('B22').Value = 'ZBToken'
$ws.Range('C22').Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
Set-TextValue $ws.Range('D22') '0.2552'
Set-TextValue $ws.Range('E22') '-0.98%'
$ws.Range('B23').Value = 'CoinExToken'
$ws.Range('C23').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
Set-TextValue $ws.Range('D23') '0.04146'
Set-TextValue $ws.Range('E23') '1.51%'
$ws.Range('B24').Value = 'BitKan'
$ws.Range('C24').Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
Set-TextValue $ws.Range('D24') '0.001302'
Set-TextValue $ws.Range('E24') '-0.33%'
Set-TextValue $ws.Range('D25') '0.0001300'
Set-TextValue $ws.Range('E25') '4.09%'
Set-TextValue $ws.Range('E26') '-20.15%'
Set-TextValue $ws.Range('D38') '0.02763'
Set-TextValue $ws.Range('E38') '6.68%'
Set-TextValue $ws.Range('D39') '0.05734'
Set-TextValue $ws.Range('E39') '0.49%'
Set-TextValue $ws.Range('D40') '0.007837'
Set-TextValue $ws.Range('E40') '3.64%'
Set-TextValue $ws.Range('D41') '0.1434'
Set-TextValue $ws.Range('E41') '2.38%'
Set-TextValue $ws.Range('D42') '0.007272'
Set-TextValue $ws.Range('E42') '-0.89%'
Set-TextValue $ws.Range('D43') '0.002109'
Set-TextValue $ws.Range('E43') '5.16%'
Set-TextValue $ws.Range('D44') '0.007913'
Set-TextValue $ws.Range('E44') '-12.73%'
Set-TextValue $ws.Range('D45') '0.3410'
Set-TextValue $ws.Range('D46') '0.00007025'
Set-TextValue $ws.Range('E46') '0.05%'
Set-TextValue $ws.Range('E47') '0.32%'
Set-TextValue $ws.Range('E48') '0.29%'
Set-TextValue $ws.Range('D49') '0.003451'
Set-TextValue $ws.Range('E49') '-1.43%'
Set-TextValue $ws.Range('D50') '0.003538'
Set-TextValue $ws.Range('E50') '1.13%'
Set-TextValue $ws.Range('D51') '0.00002105'
Set-TextValue $ws.Range('E51') '0.32%'
